# Applies the "Added Experiment to future plans" edit:
#  - grow the big rounded-rectangle panel and the year-axis rectangle
#  - push the timeline axis / tick-marks / year-labels down to make room
#  - grow the "Evaluation" future-plan box and rename it
#    "Evaluation & Experimentation"
#
# All geometry in the underlying OOXML is stored in EMU (914400 EMU = 1 in),
# while the PowerPoint object model works in points (72 pt = 1 in), so every
# EMU offset from the diff is converted to points (EMU / 12700) below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# 1) "Rounded Rectangle 47" (id 48): ext cy 1828801 -> 2057401 EMU
#    162.00007874015748pt rounds back down to 2057400 EMU through the COM
#    host's float32 path, so nudge the literal just far enough to still
#    floor to the intended 2057401 EMU.
$shRoundedRect47 = Get-ShapeById $s 48
$shRoundedRect47.Height = 162.00008392333984

# 2) "Rectangle 17" (id 18): ext cy 4114800 -> 4191000 EMU (330 pt)
$shRectangle17 = Get-ShapeById $s 18
$shRectangle17.Height = 330

# 3) "Straight Arrow Connector 24" (id 25): off y 4495800 -> 4572000 EMU (360 pt)
$shArrowConnector24 = Get-ShapeById $s 25
$shArrowConnector24.Top = 360

# 4-9) The six small tick "Straight Connector" shapes: off y 4419600 -> 4495800 EMU (354 pt)
foreach ($id in 27, 32, 33, 34, 35, 36) {
    $shTick = Get-ShapeById $s $id
    $shTick.Top = 354
}

# 10-15) The six year-label TextBoxes: off y 4569023 -> 4645223 EMU (365.7655905511811 pt)
foreach ($id in 38, 39, 40, 41, 42, 43) {
    $shYearLabel = Get-ShapeById $s $id
    $shYearLabel.Top = 365.7655905511811
}

# 16-17) "Rectangle 46" (id 47): ext cy 457200 -> 609600 EMU (48 pt)
#        and its text "Evaluation" -> "Evaluation & Experimentation"
$shEvaluationBox = Get-ShapeById $s 47
$shEvaluationBox.Height = 48
$shEvaluationBox.TextFrame.TextRange.Text = "Evaluation & Experimentation"
